# Create MeasureView, Create Quantity/Pages/Measures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 4")
$ws.Activate()

# Row 12 (entry #6): fill in date, start/stop time, delta minutes, activity,
# comment, and mark the "C" column with an "x".
$ws.Range("B12").Value = 43883
$ws.Range("C12").Value = 0.62291666666666667
$ws.Range("D12").Value = 0.75347222222222221
$ws.Range("F12").Value = 188
$ws.Range("G12").Value = "Kodutöö 4"
$ws.Range("H12").Value = "vaatasin 2. video, p. 2 tehtud"
$ws.Range("J12").Value = "x"

# Recalculate so the F18 = SUM(F7:F17) total picks up the new F12 value.
$excel.Calculate()

# Match the author's final selection after entering the data.
$ws.Range("H16").Select()
